# Add the "no two consecutive timeslots in a day" constraint support by
# updating the invigilators' time_slot_availability strings, and fix the
# selection left over from the last edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Abigail Smith (row 2): extend availability + flip lead flag to 1
$ws.Range("C2").Value = "1,2,3,4,6,7,9,10,12,13,15,16,18"
$ws.Range("D2").Value = 1

# Alice Johnson (row 3): extend availability
$ws.Range("C3").Value = "1,2,4,5,7,8,10,11,13,14"

# Annabelle Harris (row 6): extend availability
$ws.Range("C6").Value = "2,4,7,8,10,11,13,15,16,17,18"

# Andrew White (row 5): extend availability
$ws.Range("C5").Value = "1,2,3,6,7,8,9,10,12,15,16,17,18"

# Restore the active selection to C5 (single cell) as left by the author
$ws.Activate()
$ws.Range("C5").Select()
